$wb = $excel.ActiveWorkbook
$fGroups = $wb.Worksheets.Item("fGroups")

# update selection on fGroups sheet per diff
[void]$fGroups.Range("B1:G1").Select()

$ws = $wb.Worksheets.Add($null, $fGroups)
$ws.Name = "mslists"

$ws.Range("B1").Value = "as-is"
$ws.Range("C1").Value = "almost as-is"
$ws.Range("D1").Value = "implement"
$ws.Range("E1").Value = "not supported"
$ws.Range("F1").Value = "ionize"
$ws.Range("G1").Value = "done"

$ws.Range("A2").Value = "$"
$ws.Range("B2").Value = "X"

$ws.Range("A3").Value = "["
$ws.Range("C3").Value = "X"

$ws.Range("A4").Value = "[["
$ws.Range("C4").Value = "X"

$ws.Range("A5").Value = "analyses"
$ws.Range("B5").Value = "X"

$ws.Range("A6").Value = "as.data.table"
$ws.Range("C6").Value = "X"

$ws.Range("A7").Value = "averagedPeakLists"
$ws.Range("C7").Value = "X"

$ws.Range("A8").Value = "compoundViewer"
$ws.Range("E8").Value = "X"

$ws.Range("A9").Value = "filter"
$ws.Range("C9").Value = "X"

$ws.Range("A10").Value = "groupNames"
$ws.Range("B10").Value = "X"

$ws.Range("A11").Value = "initialize"
$ws.Range("C11").Value = "X"

$ws.Range("A12").Value = "length"
$ws.Range("B12").Value = "X"

$ws.Range("A13").Value = "peakLists"
$ws.Range("C13").Value = "X"

$ws.Range("A14").Value = "plotSpec"
$ws.Range("B14").Value = "X?"

$ws.Range("A15").Value = "show"
$ws.Range("C15").Value = "X"

[void]$ws.Range("D15").Select()

Write-Host "done"
